# Add a new "2023" column (T) to the table, mirroring the formatting of the
# existing "2022" column (S) for each row, and populate the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy S4:S14 (values + formatting) into T4:T14 so the new column inherits
# the same per-row cell styles as column S (header style, first-row style,
# body style, last-row style).
$ws.Range("S4:S14").Copy($ws.Range("T4"))

# Now overwrite the copied values with the actual 2023 figures.
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 99.5
$ws.Range("T6").Value = 99.426175237254469
$ws.Range("T7").Value = 99.458151211935132
$ws.Range("T8").Value = 99.44178628389156
$ws.Range("T9").Value = 99.453125
$ws.Range("T10").Value = 99.487295483676391
$ws.Range("T11").Value = 99.743589743589752
$ws.Range("T12").Value = 99.190647482014398
$ws.Range("T13").Value = 99.483321247280642
$ws.Range("T14").Value = 99.771121504627331

# Settle the cursor back on A1, matching the saved view of the final sheet.
$ws.Range("A1").Select()
